$p = $ppt.ActivePresentation

# -------------------------------------------------------------------
# Slide 2 ("Course Overview") - replace the content placeholder text
# -------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2).TextFrame.TextRange

$body2.Text = "In this course, you will learn how to:"

$p2 = $body2.InsertAfter("`rDefine DevOps")
$p3a = $body2.InsertAfter("`rIndentify")
$p3b = $body2.InsertAfter(" the challenges of using databases that are separate from other software languages and platforms")
$p4 = $body2.InsertAfter("`rInclude your database code alongside other application code in a version control system (VCS).")
$p5 = $body2.InsertAfter("`rSetup a Continuous Integration (CI) platform for your database code.")
$p6 = $body2.InsertAfter("`rWrite and include automated unit tests for your database code.")
$p7 = $body2.InsertAfter("`r- Develop an automated release process that deploys database changes to both on premise and cloud databases.")

$paras2 = $s2.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2, 6)
$paras2.IndentLevel = 2

# -------------------------------------------------------------------
# Slide 3 - "Course Structure"
# -------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Course Structure"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "This course consists of"

# -------------------------------------------------------------------
# Slide 4 - "What is DevOps?" (Donovan Brown quote + link)
# -------------------------------------------------------------------
$s4 = $p.Slides.Add(4, 2)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "What is DevOps?"

$body4 = $s4.Shapes.Item(2).TextFrame.TextRange
$body4.Text = "“DevOps is the union of people, process, and products to enable continuous delivery of value to our end users.”"
$body4.ParagraphFormat.Bullet.Type = 0

$b4p2 = $body4.InsertAfter("`rDonovan Brown")
$b4p3 = $body4.InsertAfter("`rMicrosoft Principal DevOps PM")
$b4p4a = $body4.InsertAfter("`rFrom ")
$b4p4b = $body4.InsertAfter("http://donovanbrown.com/post/what-is-devops")

$tf4 = $s4.Shapes.Item(2).TextFrame.TextRange

$line2 = $tf4.Paragraphs(2, 1)
$line2.ParagraphFormat.Alignment = 3
$line2.ParagraphFormat.Bullet.Type = 1
$line2.ParagraphFormat.Bullet.UseTextFont = 1
$line2.ParagraphFormat.Bullet.Character = 45
$line2.Font.Italic = 1

$line3 = $tf4.Paragraphs(3, 1)
$line3.ParagraphFormat.Alignment = 3
$line3.ParagraphFormat.Bullet.Type = 0
$line3.Font.Italic = 1

$line4 = $tf4.Paragraphs(4, 1)
$line4.ParagraphFormat.Alignment = 3
$line4.ParagraphFormat.Bullet.Type = 0
$line4.Font.Size = 18

$link4 = $tf4.Paragraphs(4, 1).Characters(6, 44)
$link4.ActionSettings.Item(1).Hyperlink.Address = "http://donovanbrown.com/post/what-is-devops"

# -------------------------------------------------------------------
# Slide 5 - "What is DevOps?" (longer quote)
# -------------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "What is DevOps?"

$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = "“It is very important to realize that DevOps is not a product.  You cannot buy DevOps and install it.  DevOps is not just automation or infrastructure as code.  DevOps is people following a process enabled by products to deliver value to our end users.”"
$body5.ParagraphFormat.Bullet.Type = 0
$body5.Font.Italic = 1

$b5p2 = $body5.InsertAfter("`r- Donovan Brown")
$b5p3 = $body5.InsertAfter("`r")

$tf5 = $s5.Shapes.Item(2).TextFrame.TextRange
$line5_2 = $tf5.Paragraphs(2, 1)
$line5_2.ParagraphFormat.Alignment = 3
$line5_2.ParagraphFormat.Bullet.Type = 0

# -------------------------------------------------------------------
# Slide 6 - "The Three Ways"
# -------------------------------------------------------------------
$s6 = $p.Slides.Add(6, 2)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "The Three Ways"

$body6 = $s6.Shapes.Item(2).TextFrame.TextRange
$body6.Text = "Gene Kim of The Phoenix Project and The DevOps Handbook, lists three core principles of DevOps"
$b6p2 = $body6.InsertAfter("`rSystems Thinking")
$b6p3 = $body6.InsertAfter("`rFeedback Loops")
$b6p4 = $body6.InsertAfter("`rCulture of Continuous Learning and Experimentation")
$b6p5 = $body6.InsertAfter("`r")

$tf6 = $s6.Shapes.Item(2).TextFrame.TextRange
$tf6.Paragraphs(2, 3).IndentLevel = 2

# -------------------------------------------------------------------
# Slide 7 - "Challenges for Database DevOps"
# -------------------------------------------------------------------
$s7 = $p.Slides.Add(7, 2)
$title7 = $s7.Shapes.Item(1).TextFrame.TextRange
$title7.Text = "Challenges for "
$title7.InsertAfter("Database DevOps")
